$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the current column C (ExpPoints),
# shifting ExpPoints from C to G.
$ws.Range("C1:F1").EntireColumn.Insert()

# Header row: new column headers
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"

# Copy the header style (s="1") from B1 to the new header cells
$ws.Range("B1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Keep the body cells in the new columns (rows 2-19) present but blank,
# matching the empty placeholder cells introduced by the source edit.
$ws.Range("C2:F19").Style = "Normal"
